$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Gender"
$ws.Range("E1").Value = "Stat"
